# Update countries & provincias Spain
#
# Data refresh: Rusia's case counts grew, pushing it above Reino Unido and
# Italia in the (descending, by "Casos totales") ranking; those two rows
# keep their own figures but shift down one rank. Same pattern for Ucrania
# jumping above Rumania. A handful of other countries further down the
# table simply got refreshed figures without any re-ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# --- Rusia overtakes Reino Unido and Italia (rows 6-8) ---------------------
$ws.Range("A6").Value = "Rusia"
Set-Row 6 221344 11656 39801 179534 2300 94 2009

$ws.Range("A7").Value = "Reino Unido"
Set-Row 7 219183 0 0 186984 1559 0 31855

$ws.Range("A8").Value = "Italia"
Set-Row 8 219070 0 105186 83324 1027 0 30560

# --- Singapur refreshed figures (row 28) ------------------------------------
$ws.Range("B28").Value = 23822
$ws.Range("C28").Value = 486
$ws.Range("E28").Value = 21081

# --- Ucrania overtakes Rumania (rows 37-38) ---------------------------------
$ws.Range("A37").Value = "Ucrania"
Set-Row 37 15648 416 3288 11952 207 17 408

$ws.Range("A38").Value = "Rumania"
Set-Row 38 15362 0 7051 7350 242 0 961

# --- Chequia refreshed figures (row 51) -------------------------------------
$ws.Range("D51").Value = 4482
$ws.Range("E51").Value = 3361

# --- Oman refreshed figures (row 66) ----------------------------------------
$ws.Range("B66").Value = 3573
$ws.Range("C66").Value = 174
$ws.Range("D66").Value = 1211
$ws.Range("E66").Value = 2345

# --- Estonia refreshed figures (row 83) -------------------------------------
$ws.Range("B83").Value = 1741
$ws.Range("C83").Value = 2
$ws.Range("D83").Value = 751
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 61

# --- Lituania refreshed figures (row 88) ------------------------------------
$ws.Range("B88").Value = 1485
$ws.Range("C88").Value = 6
$ws.Range("D88").Value = 833
$ws.Range("E88").Value = 602

# --- Eslovaquia refreshed figures (row 90) ----------------------------------
$ws.Range("D90").Value = 959
$ws.Range("E90").Value = 472

# --- Sri Lanka refreshed figures (row 104) ----------------------------------
$ws.Range("D104").Value = 343
$ws.Range("E104").Value = 511
